# Roadmap - CSS deck update:
#  1. Bump the cached "Today" date field (datetimeFigureOut) on the
#     slide master and every slide layout from 11-11-2025 to 12-11-2025.
#  2. Mark "Chapter 19 / Chapter 20" (row 8 of the table on slide 3) as
#     DONE, enlarging the status cell text to match the other rows.

$p = $ppt.ActivePresentation

$oldDate = "11-11-2025"
$newDate = "12-11-2025"

# --- 1. Date placeholder on the slide master ---------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 1b. Date placeholder on every slide layout -------------------------
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Slide 3 roadmap table: Chapter 19 / Chapter 20 row -> DONE -----
$slide3 = $p.Slides.Item(3)
$tableShape = $slide3.Shapes.Item(1)
$tbl = $tableShape.Table

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $contentText = $tbl.Cell($r, 1).Shape.TextFrame.TextRange.Text
    if ($contentText -match "Chapter 19" -and $contentText -match "Chapter 20") {
        $statusRange = $tbl.Cell($r, 3).Shape.TextFrame.TextRange
        $statusRange.Text = "DONE "
        $statusRange.Font.Size = 14
    }
}
